$d = $word.ActiveDocument

# Replace the whole greeting text ("Buenos días", split across multiple
# runs) with "Adios" as a single run.
$d.Content.Find.Execute("Buenos días", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Adios", 2)
